$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell text content (sharedStrings) across rows 2, 6, 9-13 ---

# Row 2
$ws.Range("C2").Value = "The fifo has a parameter 'DEPTH' that decides the depth of the fifo"
$ws.Range("D2").Value = "1. Instantiate multiple FIFOs each with a different depth`n2. Observe and check that each FIFO's depth is equal to its depth parameter"
$ws.Range("E2").Value = "DEPTH = 2, 4, 8, 16, 32 ,64"
$ws.Range("F2").Value = "FIFO's depth is equal to the depth parameter"
$ws.Rows.Item(2).RowHeight = 57.6

# Row 6
$ws.Range("F6").Value = "o_full' = 1 when fifo is full"

# Row 9
$ws.Range("D9").Value = "1. Set 'i_reset' to 1`n2. Set 'i_w_en' to 0`n3. Set 'i_r_en' to 0`n4. Set 'i_reset' to 0 after two clock cycles`n5. Check that the write pointer doen'st get incremented on the rising edge of the clock for two clock cycles`n6. set 'i_w_en' to 1`n7. Check that the write pointer gets incremented on the rising edge of the clock`n8. wait until 'o_full' is set to 1`n9. Check that the write pointer doesn't get incremented on rising edge of the clock"
$ws.Range("F9").Value = "Write pointer is incremented on rising edge of clock when fifo is not full and 'i_w_en' is set to 1"
$ws.Rows.Item(9).RowHeight = 187.2

# Row 10
$ws.Range("C10").Value = "o_fifo_w_data' is saved in the fifo cell where the write pointer is pointing to at rising edge of clock if 'i_w_en' is set to 1 and fifo is not full"
$ws.Range("D10").Value = "1. Set 'i_reset' to 1`n2. Set 'i_w_en' to 1`n3. Set 'i_fifo_w_data' to 128'h5 (value chosen arbitrarily)`n4. Set 'i_r_en' to 0`n5. Set 'i_reset' to 0 after two clock cycles`n6. Check that fifo cell where the write pointer is pointing to is equal to 128'h5 at rising edge of clock"
$ws.Range("E10").Value = "i_reset' = 0`n'i_w_en' = 1`n'i_r_en' = 0`n'i_fifo_w_data' = 128'h5"
$ws.Range("F10").Value = "FIFO cell where write pointer is pointing to is equal to 128'h5 at rising edge of the clock when 'i_w_en' is set to 1 and fifo is not full"
$ws.Rows.Item(10).RowHeight = 115.2

# Row 11
$ws.Range("C11").Value = "Read pointer gets incremented on rising edge of clock if 'i_r_en' is set to 1 and fifo is not empty"
$ws.Range("D11").Value = "1. Set 'i_w_en' to 0`n2. Set 'i_r_en' to 0`n3. Check that the read pointer doesn't get incremented on the rising edge of the clock for two clock cycles`nSet 'i_r_en' to 1`n4. Check that the read pointer gets incremented on the rising edge of the clock`n5. wait until 'o_empty' is set to 1`n6. Check that the read pointer doesn't get incremented on rising edge of the clock"
$ws.Range("F11").Value = "Read pointer is incremented on rising edge of clock when fifo is not empty and 'i_r_en' is set to 1"
$ws.Rows.Item(11).RowHeight = 158.4

# Row 12
$ws.Range("C12").Value = "o_fifo_r_data' is set to the fifo cell content where the read pointer is pointing to"
$ws.Range("D12").Value = "1. Set 'i_w_en' to 1`n2. Set 'i_fifo_w_data' to 128'h8 (value chosen arbitrarily)`n3. Wait for one clock cycle`n4. Set 'i_r_en' to 1`n5. Check that 'o_fifo_r_data' is always equal to the fifo cell where the read pointer is pointing to (in this test: 128'h8)"
$ws.Range("E12").Value = "i_reset' = 0`n'i_w_en' = 1`n'i_r_en' = 1`n'i_fifo_w_data' = 128'h8"
$ws.Range("F12").Value = "o_fifo_r_data' is always equal to the fifo cell where the read pointer is pointing to (in this test: 128'h8)"
$ws.Rows.Item(12).RowHeight = 115.2

# Row 13
$ws.Range("C13").Value = "All fifo cells are set to zeros when 'i_reset' is set to 1"
$ws.Range("D13").Value = "1. Set 'i_reset' to 1`n2. Check that all fifo cells contain zeros"
$ws.Range("F13").Value = "All fifo cells are equal to zeros"

# --- Restore quote-prefix style (s="6") on cells where Excel auto-cleared it ---
# Source cells below keep style 6 untouched by this edit; we copy their format only.
$ws.Range("E4").Copy()
$ws.Range("F6").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("C10").PasteSpecial(-4122)
$ws.Range("E5").Copy()
$ws.Range("E10").PasteSpecial(-4122)
$ws.Range("C5").Copy()
$ws.Range("C12").PasteSpecial(-4122)
$ws.Range("F4").Copy()
$ws.Range("E12").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("F12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Restore active selection ---
$ws.Range("E11").Select()
